$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HP")

# Copy formatting from row 2 (existing HP1 "geothermal heat pump" entry)
# down to the new row 5, then overwrite the values for the new HP3 unit.
$ws.Range("A2:O2").Copy()
$ws.Range("A5:O5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A5").Value = "geothermal heat pump"
$ws.Range("B5").Value = "HP3"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 10000000000
$ws.Range("E5").Value = "W"
$ws.Range("F5").Value = "USD-2015"
$ws.Range("G5").Value = 0
$ws.Range("H5").Formula = "=70.8/0.902"
$ws.Range("I5").Value = 0.49
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 20
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 5

# Make the HP sheet the active tab/selection (was Chiller before),
# matching the new activeTab + tabSelected + selection in the diff.
$ws.Activate()
$ws.Range("H6").Select()
